$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44312
$ws.Range("J2").Value = 50

# Row 3
$ws.Range("D3").Value = 44497

# Row 4
$ws.Range("D4").Value = 44498
$ws.Range("J4").Value = 40

# Row 6
$ws.Range("D6").Value = 44390
$ws.Range("J6").Value = 55
$ws.Range("K6").Value = 6000
$ws.Range("L6").Value = 6000
$ws.Range("M6").Value = 6000
$ws.Range("P6").Value = 6000

# Row 7
$ws.Range("D7").Value = 44259
$ws.Range("J7").Value = 30

# Row 8
$ws.Range("D8").Value = 44313
$ws.Range("J8").Value = 20

# Row 9
$ws.Range("D9").Value = 44176
$ws.Range("J9").Value = 10

# Row 11
$ws.Range("D11").Value = 44315
$ws.Range("J11").Value = 40

# Row 12
$ws.Range("D12").Value = 44316
$ws.Range("J12").Value = 20
$ws.Range("K12").Value = 4000
$ws.Range("L12").Value = 4000
$ws.Range("M12").Value = 4000
$ws.Range("P12").Value = 4000

# Row 13
$ws.Range("D13").Value = 44280
$ws.Range("K13").Value = 4000
$ws.Range("L13").Value = 4000
$ws.Range("M13").Value = 4000
$ws.Range("P13").Value = 4000

# Row 14
$ws.Range("D14").Value = 44301
$ws.Range("J14").Value = 40
$ws.Range("K14").Value = 3000
$ws.Range("L14").Value = 3000
$ws.Range("M14").Value = 3000
$ws.Range("P14").Value = 3000
